$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -10
    3  = -2
    5  = 0
    6  = -7
    8  = 4
    12 = -6
    16 = -3
    18 = -5
    19 = -7
    20 = -10
    21 = -7
    27 = 1
    28 = -8
    31 = -1
    33 = -7
    34 = -4
    35 = -5
    36 = 0
    37 = 1
    41 = -6
    43 = 1
    44 = -3
    48 = 0
    52 = 6
    53 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
